$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header A1: "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# 2) Column A (rows 2-14): new MaxFES fractional values replacing old Gen counts
$newA = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $newA.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $newA[$i]
}

# 3) Recompute column AZ (currently "Run 50") to become the "Mean" column,
#    averaging the 50 run columns B:AY for each data row, replacing the
#    now-removed "Run 50" run data.
for ($row = 2; $row -le 14; $row++) {
    $sum = 0.0
    for ($col = 2; $col -le 51; $col++) {
        $sum += $ws.Cells.Item($row, $col).Value()
    }
    $ws.Cells.Item($row, 52).Value = $sum / 50
}

# 4) Header AZ1: becomes "Mean" (was "Run 50", old "Mean" column BA is dropped)
$ws.Range("AZ1").Value = "Mean"

# 5) Remove the old trailing "Mean" column (BA), shifting nothing else
$ws.Range("BA1:BA14").EntireColumn.Delete()
